$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point rounding of the existing last row's timestamp (A17)
$ws.Cells.Item(17, 1).Value = 45866.75030247685

# Append the new reading as row 18
$ws.Cells.Item(18, 1).Value = 45866.79191548906
$ws.Cells.Item(18, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(18, 2).Value = 2025
$ws.Cells.Item(18, 3).Value = 31
$ws.Cells.Item(18, 4).Value = 15.28
$ws.Cells.Item(18, 5).Value = 84.25
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 1.98
$ws.Cells.Item(18, 8).Value = "E"
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = "19:00:21"
